# This edit re-orders a block of species-observation records (rows 69-81)
# on the active worksheet. Each record's identity/location data lives in
# columns A,B,D,E,F,G,H,Q,R; the remaining columns (C,I,P,S,T,U,V,W,...)
# describe the shared site/visit info and are left untouched.
#
# The new row order is a permutation of the old one, so we must first snapshot
# every original value (reading them all from the sheet as it currently
# stands) before overwriting any cell - otherwise a later read could pick up
# an already-overwritten (new) value instead of the original one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")
$rows = 69..81

# Destination row -> source row: the data that ends up in row <key> is the
# original data that used to sit in row <value>.
$rowMap = @{
    69 = 73
    70 = 74
    71 = 79
    72 = 78
    73 = 81
    74 = 70
    75 = 71
    76 = 69
    77 = 80
    78 = 77
    79 = 72
    80 = 76
    81 = 75
}

# Step 1: snapshot the original values of every relevant cell.
$orig = @{}
foreach ($r in $rows) {
    foreach ($c in $cols) {
        $addr = "$c$r"
        $orig[$addr] = $ws.Range($addr).Value2
    }
}

# Step 2: write the permuted values back, based solely on the snapshot.
foreach ($dst in $rows) {
    $src = $rowMap[$dst]
    foreach ($c in $cols) {
        $ws.Range("$c$dst").Value2 = $orig["$c$src"]
    }
}
